$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVTs")

# 1. Expand/Collapse feature checklist row (C14): extend the tooltip description
#    to mention bars and individual target points.
$ws.Range("C14").Value = "Add a tooltip for additional information on bars and individual target points(displayed only individual target)."

# 2. Data label checklist row (C11): note the new max text size constraint.
$ws.Range("C11").Value = "Update data label's color, text size, display units and decimal value(Max text size is 20)"

# 3. Update the active selection/view on the BVTs sheet to C11 (also clears the
#    previously scrolled-down top-left cell so the view resets to the top).
$ws.Activate() | Out-Null
$ws.Range("C11").Select() | Out-Null
